# Auto-generated edit script applying the diff changes
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 354
$ws.Range("F3").Value = 1257
$ws.Range("F5").Value = 2045
$ws.Range("F6").Value = 100
$ws.Range("F7").Value = 774
$ws.Range("F9").Value = 530
$ws.Range("F10").Value = 106
$ws.Range("F11").Value = 134
$ws.Range("F12").Value = 1055
$ws.Range("F13").Value = 786
$ws.Range("F14").Value = 26
$ws.Range("F15").Value = 632
$ws.Range("F16").Value = 1231
$ws.Range("F18").Value = 5
$ws.Range("F19").Value = 727
$ws.Range("F20").Value = 697
$ws.Range("F21").Value = 73
$ws.Range("F24").Value = 619
$ws.Range("F25").Value = 1189
$ws.Range("F26").Value = 113
$ws.Range("F27").Value = 175
$ws.Range("F28").Value = 4837
$ws.Range("F29").Value = 222
$ws.Range("F31").Value = 1373
$ws.Range("F32").Value = 5741
$ws.Range("F33").Value = 945
$ws.Range("F34").Value = 569
$ws.Range("F35").Value = 52
$ws.Range("F37").Value = 1035
$ws.Range("F40").Value = 633
$ws.Range("F42").Value = 23
$ws.Range("F44").Value = 5
$ws.Range("F46").Value = 86
$ws.Range("F47").Value = 289
$ws.Range("F49").Value = 13

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 8
$ws.Range("G5").Value = "不可售"
$ws.Range("F6").Value = 2071
$ws.Range("F7").Value = 69
$ws.Range("F9").Value = 108
$ws.Range("F12").Value = 92
$ws.Range("F16").Value = 636
$ws.Range("F17").Value = 636
$ws.Range("F29").Value = 1705
$ws.Range("F33").Value = 4
$ws.Range("F38").Value = 58
$ws.Range("F43").Value = 468

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 756
$ws.Range("F6").Value = 725
$ws.Range("F7").Value = 343
$ws.Range("F8").Value = 191

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 756
$ws.Range("F3").Value = 354
$ws.Range("F6").Value = 1257
$ws.Range("F7").Value = 725
$ws.Range("F8").Value = 725
$ws.Range("F10").Value = 343
$ws.Range("F11").Value = 191
$ws.Range("F12").Value = 191
$ws.Range("F13").Value = 2071
$ws.Range("F14").Value = 2045
$ws.Range("F15").Value = 69
$ws.Range("F16").Value = 774
$ws.Range("F17").Value = 108
$ws.Range("F19").Value = 106
$ws.Range("F20").Value = 134
$ws.Range("F21").Value = 1055
$ws.Range("F22").Value = 786
$ws.Range("F23").Value = 26
$ws.Range("F24").Value = 1231
$ws.Range("F25").Value = 92
$ws.Range("F27").Value = 727
$ws.Range("F29").Value = 697
$ws.Range("F30").Value = 636
$ws.Range("F31").Value = 619
$ws.Range("F32").Value = 113
$ws.Range("F35").Value = 175
$ws.Range("F37").Value = 4837
$ws.Range("F38").Value = 1373
$ws.Range("F39").Value = 5741
$ws.Range("F40").Value = 945
$ws.Range("F41").Value = 1705
$ws.Range("F42").Value = 569
$ws.Range("F43").Value = 52
$ws.Range("F44").Value = 1035
$ws.Range("F45").Value = 633
$ws.Range("F46").Value = 58
$ws.Range("F47").Value = 23
$ws.Range("F50").Value = 468
$ws.Range("F51").Value = 86
